$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '29.074.14'
$ws.Cells.Item(2, 5).Value = '  +0.67%  '

$ws.Cells.Item(3, 4).Value = '1.834.17'
$ws.Cells.Item(3, 5).Value = '  +0.71%  '

$ws.Cells.Item(4, 4).Value = '''0.9996'
$ws.Cells.Item(4, 5).Value = '  +0.62%  '

$ws.Cells.Item(5, 4).Value = '''242.06'
$ws.Cells.Item(5, 5).Value = '  -0.31%  '

$ws.Cells.Item(6, 4).Value = '''0.6190'
$ws.Cells.Item(6, 5).Value = '  -1.51%  '

$ws.Cells.Item(7, 4).Value = '''1.001'
$ws.Cells.Item(7, 5).Value = '  +0.67%  '

$ws.Cells.Item(8, 4).Value = '''0.07452'
$ws.Cells.Item(8, 5).Value = '  +0.10%  '

$ws.Cells.Item(9, 4).Value = '''0.2925'
$ws.Cells.Item(9, 5).Value = '  -0.10%  '

$ws.Cells.Item(10, 4).Value = '''23.10'
$ws.Cells.Item(10, 5).Value = '  +0.53%  '

$ws.Cells.Item(11, 4).Value = '''0.07670'
$ws.Cells.Item(11, 5).Value = '  +0.03%  '

$ws.Cells.Item(12, 4).Value = '1.834.00'
$ws.Cells.Item(12, 5).Value = '  +0.59%  '

$ws.Cells.Item(13, 4).Value = '''5.010'
$ws.Cells.Item(13, 5).Value = '  +0.79%  '

$ws.Cells.Item(14, 4).Value = '''0.6742'
$ws.Cells.Item(14, 5).Value = '  +1.37%  '

$ws.Cells.Item(15, 4).Value = '''83.02'
$ws.Cells.Item(15, 5).Value = '  +0.32%  '

$ws.Cells.Item(16, 4).Value = '''0.000009146'
$ws.Cells.Item(16, 5).Value = '  -5.34%  '

$ws.Cells.Item(17, 4).Value = '''5.914'
$ws.Cells.Item(17, 5).Value = '  -1.51%  '

$ws.Cells.Item(18, 4).Value = '29.058.77'
$ws.Cells.Item(18, 5).Value = '  +0.49%  '

$ws.Cells.Item(19, 4).Value = '2.088.53'
$ws.Cells.Item(19, 5).Value = '  +0.73%  '

$ws.Cells.Item(20, 4).Value = '''241.15'
$ws.Cells.Item(20, 5).Value = '  +7.24%  '

$ws.Cells.Item(21, 4).Value = '''12.72'
$ws.Cells.Item(21, 5).Value = '  +1.63%  '

$ws.Cells.Item(22, 5).Value = '  +0.79%  '

$ws.Cells.Item(23, 4).Value = '''7.212'
$ws.Cells.Item(23, 5).Value = '  +1.53%  '

$ws.Cells.Item(24, 4).Value = '''1.002'
$ws.Cells.Item(24, 5).Value = '  +0.69%  '

$ws.Cells.Item(25, 4).Value = '''158.97'
$ws.Cells.Item(25, 5).Value = '  -0.54%  '

$ws.Cells.Item(26, 4).Value = '''0.1412'
$ws.Cells.Item(26, 5).Value = '  +0.30%  '

$ws.Cells.Item(27, 4).Value = '''8.510'
$ws.Cells.Item(27, 5).Value = '  +0.38%  '

$ws.Cells.Item(28, 5).Value = '  +0.22%  '

$ws.Cells.Item(29, 4).Value = '''1.501'
$ws.Cells.Item(29, 5).Value = '  +0.29%  '

$ws.Cells.Item(30, 4).Value = '''0.05623'

$ws.Cells.Item(31, 2).Value = 'Filecoin'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(31, 4).Value = '''4.135'
$ws.Cells.Item(31, 5).Value = '  +0.78%  '

$ws.Cells.Item(32, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(32, 4).Value = '''4.117'
$ws.Cells.Item(32, 5).Value = '  +1.91%  '

$ws.Cells.Item(33, 4).Value = '''1.203'
$ws.Cells.Item(33, 5).Value = '  +0.65%  '

$ws.Cells.Item(34, 4).Value = '''1.843'
$ws.Cells.Item(34, 5).Value = '  -0.28%  '

$ws.Cells.Item(35, 4).Value = '''0.7418'
$ws.Cells.Item(35, 5).Value = '  -0.03%  '

$ws.Cells.Item(36, 5).Value = '  +1.06%  '

$ws.Cells.Item(37, 4).Value = '''2.658'
$ws.Cells.Item(37, 5).Value = '  +2.12%  '

$ws.Cells.Item(38, 4).Value = '''2.774'
$ws.Cells.Item(38, 5).Value = '  +1.34%  '

$ws.Cells.Item(39, 4).Value = '''0.01786'
$ws.Cells.Item(39, 5).Value = '  +0.75%  '

$ws.Cells.Item(40, 4).Value = '1.210.24'
$ws.Cells.Item(40, 5).Value = '  -2.08%  '

$ws.Cells.Item(41, 4).Value = '''6.402'
$ws.Cells.Item(41, 5).Value = '  -3.93%  '

$ws.Cells.Item(42, 4).Value = '''0.8997'
$ws.Cells.Item(42, 5).Value = '  +0.47%  '

$ws.Cells.Item(43, 4).Value = '''0.9998'
$ws.Cells.Item(43, 5).Value = '  +0.56%  '

$ws.Cells.Item(44, 4).Value = '''101.53'
$ws.Cells.Item(44, 5).Value = '  +0.35%  '

$ws.Cells.Item(45, 4).Value = '1.984.20'
$ws.Cells.Item(45, 5).Value = '  +0.64%  '

$ws.Cells.Item(46, 4).Value = '''65.48'
$ws.Cells.Item(46, 5).Value = '  +1.07%  '

$ws.Cells.Item(47, 4).Value = '''0.5092'
$ws.Cells.Item(47, 5).Value = '  +0.66%  '

$ws.Cells.Item(48, 4).Value = '''0.4064'
$ws.Cells.Item(48, 5).Value = '  +0.74%  '

$ws.Cells.Item(49, 5).Value = '  -4.62%  '

$ws.Cells.Item(50, 4).Value = '''9.123'
$ws.Cells.Item(50, 5).Value = '  +2.32%  '

$ws.Cells.Item(51, 4).Value = '''0.05815'
$ws.Cells.Item(51, 5).Value = '  +0.70%  '

$wb.Save()